$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Reverse the order of the period labels in E16:E22 (text values),
# and swap the mora values in F16/F22 to match.
$ws.Range("E16").Value = "2211"
$ws.Range("E17").Value = "2210"
$ws.Range("E18").Value = "2209"
$ws.Range("E19").Value = "2208"
$ws.Range("E20").Value = "2207"
$ws.Range("E21").Value = "2206"
$ws.Range("E22").Value = "2205"

$ws.Range("F16").Value = 28000
$ws.Range("F22").Value = 40000
